# Add a blue (accent1) highlight to the two "problem area" sub-bullets on
# the "Personnel" slide (slide 12, "Content Placeholder 2" shape):
#   - "Limited availability Feb-Mar 2019"
#   - "New developer – needs time to familiarize himself (ready ~Mar 2019)"
#
# Both runs already carry i="1" (italic); we only need to add a solidFill
# using the accent1 theme color (msoThemeColorAccent1 = 5).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

$targets = @(
    "Limited availability Feb-Mar 2019",
    "New developer – needs time to familiarize himself (ready ~Mar 2019)"
)

$count = $tr.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $text = $para.Text.TrimEnd()
    if ($targets -contains $text) {
        $para.Font.Color.ObjectThemeColor = 5  # msoThemeColorAccent1
    }
}
